# Scheduled runner update: refresh Leve profit calculations across all job sheets.
# Mirrors the upstream "Zodiark_Profits" data refresh - only the H/I/J/K/L (price)
# and M/N (profit) columns move; item identity columns A-G are untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null

$ws.Range("H93").Value = 19199.334
$ws.Range("J93").Value = 19199.334
$ws.Range("L93").Value = 19199.334
$ws.Range("N93").Value = -24191.334

$ws.Range("H116").Value = 3639.7778
$ws.Range("I116").Value = 3336
$ws.Range("K116").Value = 3336
$ws.Range("M116").Value = 106

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H76").Value = 72637.42999999999
$ws.Range("J76").Value = 70593.836
$ws.Range("L76").Value = 70593.836
$ws.Range("N76").Value = -71269.836

$ws.Range("H79").Value = 72637.42999999999
$ws.Range("J79").Value = 70593.836
$ws.Range("L79").Value = 70593.836
$ws.Range("N79").Value = -72933.836

$ws.Range("H92").Value = 98387.25
$ws.Range("J92").Value = 98387.25
$ws.Range("L92").Value = 98387.25
$ws.Range("N92").Value = -103379.25

$ws.Range("H97").Value = 341.89474
$ws.Range("I97").Value = 312.93332
$ws.Range("J97").Value = 450.5
$ws.Range("K97").Value = 312.93332
$ws.Range("L97").Value = 450.5
$ws.Range("M97").Value = 183.06668
$ws.Range("N97").Value = -1442.5

$ws.Range("H107").Value = 36867
$ws.Range("J107").Value = 36867
$ws.Range("L107").Value = 36867
$ws.Range("N107").Value = -44547

$ws.Range("H112").Value = 16380.333
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 16380.333
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 16380.333
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -19334.333

$ws.Range("H122").Value = 5510.1333
$ws.Range("I122").Value = 5395.59
$ws.Range("J122").Value = 6254.6665
$ws.Range("K122").Value = 16186.77
$ws.Range("L122").Value = 18763.9995
$ws.Range("M122").Value = -13736.77
$ws.Range("N122").Value = -23663.9995

$ws.Range("H124").Value = 23964.5
$ws.Range("J124").Value = 23964.5
$ws.Range("L124").Value = 23964.5
$ws.Range("N124").Value = -33784.5

$ws.Range("H135").Value = 85300
$ws.Range("J135").Value = 85300
$ws.Range("L135").Value = 85300
$ws.Range("N135").Value = -95440

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H81").Value = 42095
$ws.Range("J81").Value = 42095
$ws.Range("L81").Value = 42095
$ws.Range("N81").Value = -44217

$ws.Range("H84").Value = 42095
$ws.Range("J84").Value = 42095
$ws.Range("L84").Value = 126285
$ws.Range("N84").Value = -136893

$ws.Range("H100").Value = 57950
$ws.Range("J100").Value = 57950
$ws.Range("L100").Value = 57950
$ws.Range("N100").Value = -60114

$ws.Range("H110").Value = 20652
$ws.Range("J110").Value = 20652
$ws.Range("L110").Value = 20652
$ws.Range("N110").Value = -28832

$ws.Range("H112").Value = 42772.5
$ws.Range("I112").Value = 75000
$ws.Range("J112").Value = 10545
$ws.Range("K112").Value = 75000
$ws.Range("L112").Value = 10545
$ws.Range("M112").Value = -73523
$ws.Range("N112").Value = -13499

$ws.Range("H130").Value = 64798
$ws.Range("J130").Value = 64798
$ws.Range("L130").Value = 64798
$ws.Range("N130").Value = -74838

$ws.Range("H134").Value = 2538.9048
$ws.Range("I134").Value = 2015.85
$ws.Range("J134").Value = 13000
$ws.Range("K134").Value = 6047.549999999999
$ws.Range("L134").Value = 39000
$ws.Range("M134").Value = -3512.549999999999
$ws.Range("N134").Value = -44070

$ws.Range("H135").Value = 79491.664
$ws.Range("J135").Value = 79491.664
$ws.Range("L135").Value = 79491.664
$ws.Range("N135").Value = -89631.664

# ---------------------------------------------------------------------------
# CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H28").Value = 38960.25
$ws.Range("J28").Value = 38960.25
$ws.Range("L28").Value = 38960.25
$ws.Range("N28").Value = -39450.25

$ws.Range("H31").Value = 1177.4286
$ws.Range("I31").Value = 1177.4286
$ws.Range("K31").Value = 1177.4286
$ws.Range("M31").Value = -882.4286

$ws.Range("H34").Value = 1177.4286
$ws.Range("I34").Value = 1177.4286
$ws.Range("K34").Value = 1177.4286
$ws.Range("M34").Value = -975.4286

$ws.Range("H59").Value = 65138
$ws.Range("J59").Value = 65138
$ws.Range("L59").Value = 65138
$ws.Range("N59").Value = -67428

$ws.Range("H95").Value = 17230.6
$ws.Range("J95").Value = 17230.6
$ws.Range("L95").Value = 17230.6
$ws.Range("N95").Value = -22722.6

$ws.Range("H96").Value = 10381.333
$ws.Range("J96").Value = 10381.333
$ws.Range("L96").Value = 10381.333
$ws.Range("N96").Value = -15873.333

# ---------------------------------------------------------------------------
# CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H58").Value = 3624.75
$ws.Range("I58").Value = 3499.6667
$ws.Range("K58").Value = 10499.0001
$ws.Range("M58").Value = -10371.0001

$ws.Range("H120").Value = 11676.167
$ws.Range("I120").Value = 11676.167
$ws.Range("K120").Value = 35028.501
$ws.Range("M120").Value = -30190.501

$ws.Range("H131").Value = 1645.5834
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 47811.43
$ws.Range("I70").Value = 82622.28999999999
$ws.Range("J70").Value = 13000.571
$ws.Range("K70").Value = 82622.28999999999
$ws.Range("L70").Value = 13000.571
$ws.Range("M70").Value = -82352.28999999999
$ws.Range("N70").Value = -13540.571

$ws.Range("H73").Value = 47811.43
$ws.Range("I73").Value = 82622.28999999999
$ws.Range("J73").Value = 13000.571
$ws.Range("K73").Value = 82622.28999999999
$ws.Range("L73").Value = 13000.571
$ws.Range("M73").Value = -81686.28999999999
$ws.Range("N73").Value = -14872.571

$ws.Range("H132").Value = 3022.2666
$ws.Range("I132").Value = 2652.25
$ws.Range("K132").Value = 7956.75
$ws.Range("M132").Value = -5426.75

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H101").Value = 67812.664
$ws.Range("J101").Value = 67812.664
$ws.Range("L101").Value = 67812.664
$ws.Range("N101").Value = -74302.664

$ws.Range("H106").Value = 20430
$ws.Range("J106").Value = 20430
$ws.Range("L106").Value = 20430
$ws.Range("N106").Value = -22954

$ws.Range("H110").Value = 24100.5
$ws.Range("J110").Value = 24100.5
$ws.Range("L110").Value = 24100.5
$ws.Range("N110").Value = -32280.5

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H92").Value = 70000
$ws.Range("J92").Value = 70000
$ws.Range("L92").Value = 70000
$ws.Range("N92").Value = -74992

$ws.Range("H105").Value = 30213.428
$ws.Range("J105").Value = 30213.428
$ws.Range("L105").Value = 30213.428
$ws.Range("N105").Value = -37201.428

$ws.Range("H119").Value = 77530.336
$ws.Range("J119").Value = 77530.336
$ws.Range("L119").Value = 77530.336
$ws.Range("N119").Value = -87206.336
